$wb = $excel.ActiveWorkbook

# Clone the existing "Sedan_HambaLG" sheet to create the new "FSAE_Achilles" sheet,
# placed right after it, so it inherits the same layout/formatting.
$src = $wb.Worksheets.Item("Sedan_HambaLG")
$src.Copy($null, $src)

# The copy becomes active and is placed immediately after the source sheet.
$newSheet = $wb.ActiveSheet
$newSheet.Name = "FSAE_Achilles"

# Sheet-specific tweaks for the new template.
$newSheet.Range("H3").Value = "FSAE_Achilles"
$newSheet.Range("H6").Value = 0.25

# Select H4 on the new sheet (matches the source sheet's saved selection state)
$newSheet.Range("H4").Select()

# Make the new sheet the active tab.
$newSheet.Activate()
